# Applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.615.56'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '2.050.44'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.33'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.70%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '63.06'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.369'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0753'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('E12').Value = '  -3.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.928'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.20%  '
$ws.Range('D15').Value = '2.351.62'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.45'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.53%  '
$ws.Range('D17').Value = '2.044.00'
$ws.Range('E17').Value = '  -1.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.85'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.25%  '
$ws.Range('D19').Value = '36.548.45'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.93%  '
$ws.Range('D21').Value = '0.0₃0861'
$ws.Range('E21').Value = '  -3.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '237.99'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.36%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('E25').Value = '  -2.55%  '
$ws.Range('E26').Value = '  +2.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '164.96'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.122'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.21'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0605'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.45'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.95%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0872'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.36%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.22'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.11'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.10%  '
$ws.Range('E40').Value = '  -6.21%  '
$ws.Range('E41').Value = '  -1.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0217'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.28%  '
$ws.Range('E43').Value = '  -4.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '94.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.69%  '
$ws.Range('E45').Value = '  -5.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.05'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.89%  '
$ws.Range('D47').Value = '1.380.58'
$ws.Range('E47').Value = '  +4.69%  '
$ws.Range('E48').Value = '  +7.63%  '
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.57%  '
